# Weekly driver report update for 2025-04-20
# Updates the "Bad Drivers" totals (row 4/5) and re-sorts / refreshes the
# "Good Drivers" table (rows 13-18) with the latest driver vintage data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Bad Drivers section
# ---------------------------------------------------------------------
$ws.Range("B4").Value = 113
$ws.Range("C4").Value = 3322
$ws.Range("D4").Value = 98.8

$ws.Range("B5").Value = 114
$ws.Range("C5").Value = 3324

# ---------------------------------------------------------------------
# Good Drivers section (rows 13-18), re-sorted by Driver Vintage (desc)
# ---------------------------------------------------------------------
$goodDrivers = @(
    @{ Row = 13; Name = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4";  Count = 445055; Pct = 99.90000000000001; Vintage = "2024-11-10" },
    @{ Row = 14; Name = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9";   Count = 77849;  Pct = 99.90000000000001; Vintage = "2021-08-18" },
    @{ Row = 15; Name = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1";   Count = 34244;  Pct = 100;               Vintage = "2021-04-27" },
    @{ Row = 16; Name = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2";  Count = 59673;  Pct = 100;               Vintage = "2020-08-05" },
    @{ Row = 17; Name = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6";   Count = 113652; Pct = 100;               Vintage = "2020-01-06" },
    @{ Row = 18; Name = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1";   Count = 56018;  Pct = 100;               Vintage = "2019-12-14" }
)

foreach ($drv in $goodDrivers) {
    $r = $drv.Row
    $ws.Cells.Item($r, 1).Value = $drv.Name
    $ws.Cells.Item($r, 2).Value = $drv.Count
    $ws.Cells.Item($r, 4).Value = $drv.Pct

    $eCell = $ws.Cells.Item($r, 5)
    $existing = $eCell.Value2
    if ($existing -ne $drv.Vintage) {
        # Leading apostrophe forces the date-like string to stay plain text
        # instead of being auto-converted into a date serial number.
        $eCell.Value = "'" + $drv.Vintage
    }
}
